$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - "Numero Empleado" was blank, now set to 1001 (kept as text)
$ws.Range("A4").Value = "'1001"

# Rows 5-11 - "Numero Empleado" 333385 -> 1002 (kept as text)
$ws.Range("A5:A11").Value = "'1002"

# Row 12 - new employee record: Numero Empleado 1003, Empleado "Ivan Arrieta Arrieta"
$ws.Range("A12").Value = "'1003"
$ws.Range("B12").Value = "Ivan Arrieta Arrieta"
